$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row for the LeetCode "Koko Eating Bananas" (#875) entry, appended right
# after the existing last row (26).
$row = 27

# Match the formatting (styles, number formats, wrap text, etc.) of the row
# directly above it (row 26) before filling in the values.
$ws.Range("A26:I26").Copy() | Out-Null
$ws.Range("A27:I27").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($row, 1).Value = 875
$ws.Cells.Item($row, 2).Value = "Koko Eating Bananas"
$ws.Cells.Item($row, 3).Value = "#two-pointers #array #binary-search #必背 "
$ws.Cells.Item($row, 4).Value = "medium"
$ws.Cells.Item($row, 5).Value = 3
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 10
$ws.Cells.Item($row, 8).Value = 45838
$ws.Cells.Item($row, 9).Value = 45838

# The row holds a multi-line tag string, so it needs the taller wrapped
# height like the other similarly-tagged rows above.
$ws.Rows.Item($row).RowHeight = 51

# Update the active selection / view like the saved workbook.
$ws.Range("I27").Select() | Out-Null
